# Add files via upload
# Updates the CIFAR-100 global-results sheet: refreshed metric values for
# rows 3 (Fine-tuning), 5 (Rehearsal 0.1), 6 (Rehearsal 0.5), 7 (EWC) and
# 8 (LwF), plus a full box border around the merged header row (A1:L1) so
# every header cell shares one consistent bordered/filled/centered style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: Fine-tuning ---
$ws.Range("B3").Value = 54.46
$ws.Range("C3").Value = 50.28
$ws.Range("D3").Value = 55.16
$ws.Range("E3").Value = 54.58
$ws.Range("F3").Value = 55.78
$ws.Range("G3").Value = 56
$ws.Range("H3").Value = 54.16
$ws.Range("I3").Value = 54.96
$ws.Range("J3").Value = 54.2
$ws.Range("K3").Value = 55.44
$ws.Range("L3").Value = 54.7

# --- Row 4: Joint datasets (only column C changes) ---
$ws.Range("C4").Value = 50.28

# --- Row 5: Rehearsal 0.1 ---
$ws.Range("B5").Value = 27.23
$ws.Range("C5").Value = 50.28
$ws.Range("D5").Value = 27.58
$ws.Range("E5").Value = 27.29
$ws.Range("F5").Value = 27.89
$ws.Range("G5").Value = 28
$ws.Range("H5").Value = 27.08
$ws.Range("I5").Value = 27.48
$ws.Range("J5").Value = 27.1
$ws.Range("K5").Value = 27.72
$ws.Range("L5").Value = 27.35

# --- Row 6: Rehearsal 0.5 ---
$ws.Range("C6").Value = 50.28
$ws.Range("D6").Value = 13.06
$ws.Range("E6").Value = 31.26
$ws.Range("F6").Value = 40.26
$ws.Range("G6").Value = 3.54
$ws.Range("H6").Value = 40.04
$ws.Range("I6").Value = 11.6
$ws.Range("J6").Value = 44.7
$ws.Range("K6").Value = 5.2
$ws.Range("L6").Value = 21.04

# --- Row 7: EWC ---
$ws.Range("B7").Value = 58.98
$ws.Range("C7").Value = 50.28
$ws.Range("D7").Value = 57.04
$ws.Range("E7").Value = 52.88
$ws.Range("F7").Value = 53.44
$ws.Range("G7").Value = 12.54
$ws.Range("H7").Value = 40.6
$ws.Range("I7").Value = 54.7
$ws.Range("J7").Value = 33.5
$ws.Range("K7").Value = 58.32
$ws.Range("L7").Value = 48.94

# --- Row 8: LwF ---
$ws.Range("B8").Value = 29.49
$ws.Range("C8").Value = 50.28
$ws.Range("D8").Value = 35.05
$ws.Range("E8").Value = 42.07
$ws.Range("F8").Value = 46.84999999999999
$ws.Range("G8").Value = 8.039999999999999
$ws.Range("H8").Value = 40.32
$ws.Range("I8").Value = 33.15
$ws.Range("J8").Value = 39.1
$ws.Range("K8").Value = 31.76
$ws.Range("L8").Value = 34.98999999999999

# --- Header row formatting: give the whole merged title band A1:L1 a
#     matching thin box border on every side (previously only A1 carried a
#     left-only border) so the header reads as one bordered/filled bar. ---
$header = $ws.Range("A1:L1")
$header.Borders.LineStyle = 1
$header.Borders.Weight = 2
$header.Font.Bold = $true
$header.Interior.Color = $ws.Range("A1").Interior.Color
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4108
